# Generate Report for Handoff
# Updates the "b.md" row (row 3) on every sheet to reflect that a new
# handoff has been generated:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The latest handoff file / datetime columns point at the new handoff
#     package (b.63290e5768f688058c7b37413b0a5c26c308f864.<locale>.xlf)

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-22 16:36:47"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-22 16:36:43"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-22 16:36:47"
